$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose new Price text is a plain-looking number need an explicit
# Text format first, otherwise Excel would coerce the literal into a
# numeric cell instead of keeping it as the original text value.
$textFormatRows = @(4,5,6,7,8,9,10,11,13,14,15,16,18,19,20,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $textFormatRows) { $ws.Range("D$r").NumberFormat = "@" }

$ws.Range("D2").Value = "30.416.47"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.928.29"
$ws.Range("E3").Value = "  +3.89%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "241.01"
$ws.Range("E5").Value = "  +3.10%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").Value = "0.4758"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.2858"
$ws.Range("E8").Value = "  +3.81%  "
$ws.Range("D9").Value = "0.06603"
$ws.Range("E9").Value = "  +4.17%  "
$ws.Range("D10").Value = "19.16"
$ws.Range("E10").Value = "  +8.39%  "
$ws.Range("D11").Value = "105.51"
$ws.Range("E11").Value = "  +24.40%  "
$ws.Range("D12").Value = "1.920.59"
$ws.Range("E12").Value = "  +3.32%  "
$ws.Range("D13").Value = "0.07602"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").Value = "5.127"
$ws.Range("E14").Value = "  +2.86%  "
$ws.Range("D15").Value = "0.6583"
$ws.Range("E15").Value = "  +4.81%  "
$ws.Range("D16").Value = "301.68"
$ws.Range("E16").Value = "  +22.21%  "
$ws.Range("D17").Value = "30.416.94"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "12.91"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").Value = "0.000007516"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").Value = "2.168.52"
$ws.Range("E21").Value = "  +3.44%  "
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "5.255"
$ws.Range("E23").Value = "  +6.31%  "
$ws.Range("D24").Value = "6.323"
$ws.Range("E24").Value = "  +6.74%  "
$ws.Range("D25").Value = "168.46"
$ws.Range("E25").Value = "  +2.78%  "
$ws.Range("D26").Value = "9.233"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").Value = "19.74"
$ws.Range("E27").Value = "  +9.79%  "
$ws.Range("D28").Value = "2.010"
$ws.Range("E28").Value = "  +6.83%  "
$ws.Range("D29").Value = "0.1128"
$ws.Range("E29").Value = "  +9.74%  "
$ws.Range("D30").Value = "1.352"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "4.101"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("D32").Value = "3.923"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").Value = "0.05014"
$ws.Range("E33").Value = "  +3.23%  "
$ws.Range("D34").Value = "0.7422"
$ws.Range("E34").Value = "  +5.65%  "
$ws.Range("D35").Value = "1.149"
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("D36").Value = "0.9992"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "2.730"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "0.01952"
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("D39").Value = "2.699"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").Value = "2.047"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("D41").Value = "0.8723"
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("D42").Value = "107.50"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("D43").Value = "5.798"
$ws.Range("E43").Value = "  +4.51%  "
$ws.Range("D44").Value = "70.14"
$ws.Range("E44").Value = "  +11.09%  "
$ws.Range("D45").Value = "0.9999"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "0.4135"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").Value = "7.216"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").Value = "9.264"
$ws.Range("E48").Value = "  +7.64%  "
$ws.Range("D49").Value = "34.80"
$ws.Range("E49").Value = "  +3.28%  "
$ws.Range("D50").Value = "0.1204"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "0.05627"
$ws.Range("E51").Value = "  +1.79%  "
